$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample task row entirely (Math), shifting nothing else up since
# it's already the last row - then clear the remaining sample row's contents
# (English Exam), leaving the Time/Deadline formatting in place on C2:D2 so
# new tasks can be positioned in time later.
$ws.Range("A3:D3").Delete()
$ws.Range("A2:D2").ClearContents()
